$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "258.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.19%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.30%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.820"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-10.15%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.51%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.688"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.38%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8739"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.98%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9529"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.71%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1416"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.66%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03600"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07221"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.44%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03163"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.13%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.03%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001551"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.01%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.01065"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1,663.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006029"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.20%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.34%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.228"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.10%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.61%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.96%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.17%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.532"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.02%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04212"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.30%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1380"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.04%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.45%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004512"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.13%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.03%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001493"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-22.98%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03846"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.00%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005979"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.54%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.32%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002200"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.03%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01069"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.29%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005491"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.19%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.03%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1091"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "8.95%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.61%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
